$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 11
$ws1.Range("G2").Value = 45
$ws1.Range("F3").Value = 269
$ws1.Range("F4").Value = 146
$ws1.Range("F7").Value = 4921
$ws1.Range("F8").Value = 4921
$ws1.Range("F13").Value = 1118
$ws1.Range("F14").Value = 671
$ws1.Range("F15").Value = 4643
$ws1.Range("F19").Value = 231
$ws1.Range("F20").Value = 3636
$ws1.Range("F24").Value = 3384
$ws1.Range("F26").Value = 143
$ws1.Range("F28").Value = 165
$ws1.Range("F29").Value = 218
$ws1.Range("F32").Value = 84
$ws1.Range("F36").Value = 5936
$ws1.Range("F37").Value = 932
$ws1.Range("F38").Value = 445
$ws1.Range("F40").Value = 964
$ws1.Range("F42").Value = 1226
$ws1.Range("F43").Value = 129
$ws1.Range("F44").Value = 565
$ws1.Range("F46").Value = 2085
$ws1.Range("F49").Value = 736

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F7").Value = 127
$ws2.Range("F15").Value = 135
$ws2.Range("F23").Value = 773

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 269
$ws4.Range("F7").Value = 146
$ws4.Range("F10").Value = 4921
$ws4.Range("F11").Value = 4921
$ws4.Range("F16").Value = 1118
$ws4.Range("F17").Value = 671
$ws4.Range("F18").Value = 4643
$ws4.Range("F22").Value = 231
$ws4.Range("F23").Value = 3636
$ws4.Range("F24").Value = 3384
$ws4.Range("F26").Value = 143
$ws4.Range("F27").Value = 218
$ws4.Range("F30").Value = 84
$ws4.Range("F33").Value = 135
$ws4.Range("F35").Value = 5936
$ws4.Range("F36").Value = 932
$ws4.Range("F37").Value = 445
$ws4.Range("F41").Value = 964
$ws4.Range("F43").Value = 1226
$ws4.Range("F44").Value = 129
$ws4.Range("F45").Value = 565
$ws4.Range("F46").Value = 2085
$ws4.Range("F48").Value = 736
